$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows that were deleted from the export:
#   row 8  -> 004481463 / Mara     / 7555.43
#   row 3  -> 005273382 / Mvfc     / 215930.45
#   row 2  -> 004870976 / Hfr      / 499741.41
# Deleting from bottom to top keeps the remaining row numbers stable.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()

# After the deletions above, the row for 004452476 / Ivone / 502.76 is now
# row 11. Insert a new row right after it for the new account.
$ws.Rows.Item(12).Insert()

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "004979322"
$ws.Range("B12").Value = "Marilia"
$ws.Range("C12").Value = 500
